$wb = $excel.ActiveWorkbook

# Sheet ALC, row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 338.4
$ws.Range("I33").Value = 338.4
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 338.4
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -109.4
$ws.Range("N33").Value = $null

# Sheet ALC, row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 15421.143
$ws.Range("I74").Value = 7650
$ws.Range("J74").Value = 21249.5
$ws.Range("K74").Value = 7650
$ws.Range("L74").Value = 21249.5
$ws.Range("M74").Value = -6714
$ws.Range("N74").Value = -23121.5

# Sheet ALC, row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 15421.143
$ws.Range("I77").Value = 7650
$ws.Range("J77").Value = 21249.5
$ws.Range("K77").Value = 38250
$ws.Range("L77").Value = 106247.5
$ws.Range("M77").Value = -33570
$ws.Range("N77").Value = -115607.5

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2400.7646
$ws.Range("I129").Value = 1979.6666
$ws.Range("J129").Value = 2630.4546
$ws.Range("K129").Value = 5938.9998
$ws.Range("L129").Value = 7891.3638
$ws.Range("M129").Value = -938.9997999999996
$ws.Range("N129").Value = -17891.3638

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3526.2942
$ws.Range("I137").Value = 2116
$ws.Range("J137").Value = 5112.875
$ws.Range("K137").Value = 6348
$ws.Range("L137").Value = 15338.625
$ws.Range("M137").Value = -3798
$ws.Range("N137").Value = -20438.625

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2847.0952
$ws.Range("I45").Value = 2376.389
$ws.Range("K45").Value = 2376.389
$ws.Range("M45").Value = -1999.389

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2998
$ws.Range("I74").Value = 2998
$ws.Range("K74").Value = 2998
$ws.Range("M74").Value = -2124

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2998
$ws.Range("I77").Value = 2998
$ws.Range("K77").Value = 14990
$ws.Range("M77").Value = -10622

# Sheet CRP, row 2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 325.07144
$ws.Range("I2").Value = 4.6666665
$ws.Range("J2").Value = 2247.5
$ws.Range("K2").Value = 4.6666665
$ws.Range("L2").Value = 2247.5
$ws.Range("M2").Value = 108.3333335
$ws.Range("N2").Value = -2473.5

# Sheet CRP, row 98
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 50780
$ws.Range("J98").Value = 50780
$ws.Range("L98").Value = 50780
$ws.Range("N98").Value = -55272

# Sheet CUL, row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 41.086956
$ws.Range("J2").Value = 47.7
$ws.Range("L2").Value = 286.2
$ws.Range("N2").Value = -512.2

# Sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 147147.33
$ws.Range("I4").Value = 600020.2
$ws.Range("J4").Value = 5624.5625
$ws.Range("K4").Value = 1800060.6
$ws.Range("L4").Value = 16873.6875
$ws.Range("M4").Value = -1799948.6
$ws.Range("N4").Value = -17097.6875

# Sheet CUL, row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 88.75
$ws.Range("I17").Value = 88.75
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 266.25
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -97.25
$ws.Range("N17").Value = $null

# Sheet CUL, row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3132
$ws.Range("J34").Value = 4648.25
$ws.Range("L34").Value = 13944.75
$ws.Range("N34").Value = -14112.75

# Sheet CUL, row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6841.3
$ws.Range("J39").Value = 9713.714
$ws.Range("L39").Value = 29141.142
$ws.Range("N39").Value = -29729.142

# Sheet CUL, row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 553
$ws.Range("I86").Value = 350
$ws.Range("K86").Value = 1050
$ws.Range("M86").Value = 136

# Sheet CUL, row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 553
$ws.Range("I89").Value = 350
$ws.Range("K89").Value = 3150
$ws.Range("M89").Value = 2778

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3291.6667
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3291.6667
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 29625.0003
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -34685.0003

# Sheet CUL, row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 6747.75
$ws.Range("I134").Value = 1330.5
$ws.Range("K134").Value = 3991.5
$ws.Range("M134").Value = 1078.5

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2938.9092
$ws.Range("I102").Value = 2422.5557
$ws.Range("K102").Value = 2422.5557
$ws.Range("M102").Value = -800.5556999999999

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 48411.76
$ws.Range("I132").Value = 55978.617
$ws.Range("K132").Value = 167935.851
$ws.Range("M132").Value = -165405.851

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8036.5
$ws.Range("I7").Value = 6998.3335
$ws.Range("K7").Value = 6998.3335
$ws.Range("M7").Value = -6886.3335

# Sheet LTW, row 64
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25450

# Sheet LTW, row 67
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26560

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8036.5
$ws.Range("I126").Value = 6998.3335
$ws.Range("K126").Value = 20995.0005
$ws.Range("M126").Value = -18525.0005

# Sheet WVR, row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2442.4285
$ws.Range("I2").Value = 3059.8
$ws.Range("J2").Value = 899
$ws.Range("K2").Value = 3059.8
$ws.Range("L2").Value = 899
$ws.Range("M2").Value = -2947.8
$ws.Range("N2").Value = -1123

# Sheet WVR, row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 11099.4
$ws.Range("I4").Value = 16728.5
$ws.Range("J4").Value = 2655.75
$ws.Range("K4").Value = 16728.5
$ws.Range("L4").Value = 2655.75
$ws.Range("M4").Value = -16615.5
$ws.Range("N4").Value = -2881.75

# Sheet WVR, row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 76766.22
$ws.Range("J54").Value = 76766.22
$ws.Range("L54").Value = 76766.22
$ws.Range("N54").Value = -77806.22

# Sheet WVR, row 63
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 29892.428
$ws.Range("J63").Value = 29892.428
$ws.Range("L63").Value = 29892.428
$ws.Range("N63").Value = -31140.428

# Sheet WVR, row 66
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 29892.428
$ws.Range("J66").Value = 29892.428
$ws.Range("L66").Value = 89677.284
$ws.Range("N66").Value = -95917.284

# Sheet WVR, row 95
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 31857
$ws.Range("J95").Value = 31857
$ws.Range("L95").Value = 31857
$ws.Range("N95").Value = -37349

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3775.0527
$ws.Range("I126").Value = 2266.1428
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 6798.428400000001
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -4328.428400000001
$ws.Range("N126").Value = -28940

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3815.8572
$ws.Range("I136").Value = 2180.15
$ws.Range("K136").Value = 6540.450000000001
$ws.Range("M136").Value = -3990.450000000001
